$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename column B header: "condition" -> "replicate"
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "replicate"

# ---------------------------------------------------------------------------
# 2) Insert a new column before G. This shifts the old G..J (raw_fastq,
#    trimmed_fastq, index_table, qiime_otus) one column to the right,
#    becoming H..K, and leaves a blank new column G in their place.
# ---------------------------------------------------------------------------
$ws.Columns("G:G").Insert()

# New header for the inserted column.
$ws.Range("G1").Value = "publication_nb"

# Try to restore the (approximate) column width that the neighbouring,
# shifted-from column used to have, so G:H end up with matching widths.
$ws.Columns("G:G").ColumnWidth = 32.46

# The insert leaves a blank-but-styled placeholder cell in every row of the
# new column. The "reference" rows (9, 16, 24) never get a publication_nb
# value, so fully clear those three placeholders -- otherwise they would be
# saved out as empty styled cells that don't exist in the target sheet.
$ws.Range("G9").Clear()
$ws.Range("G16").Clear()
$ws.Range("G24").Clear()

# ---------------------------------------------------------------------------
# 3) Fill in the new "publication_nb" numbers for every sample row except
#    the "reference" rows (9, 16, 24) which are left blank, exactly as the
#    neighbouring raw_fastq column already skips a value-less cell there.
# ---------------------------------------------------------------------------
$publicationNb = @{
  2  = 126.181
  3  = 485.391
  4  = 5525.29
  5  = 38718.997
  6  = 128553.996
  7  = 264233.1456
  8  = 627093.3747
  10 = 632.378166
  11 = 4191.3941
  12 = 34308.17307
  13 = 149073.2134
  14 = 480008.9
  15 = 778000.228
  17 = 86.491743
  18 = 1135.96094801
  19 = 11349.70954
  20 = 102107.074926
  21 = 406445.9592
  22 = 1138053.3501
  23 = 1211091.30755
}

foreach ($row in $publicationNb.Keys) {
  $cell = $ws.Range("G" + $row)
  $cell.Value = $publicationNb[$row]
  $cell.NumberFormat = "0.00E+00"
  $cell.HorizontalAlignment = -4108  # xlCenter
}

# ---------------------------------------------------------------------------
# 4) Formatting: every populated cell in the sheet becomes horizontally
#    centred.  The (shifted) formula columns J:K additionally wrap text,
#    matching what used to be columns I:J before the insert.
#    Columns A:F and H:K are fully populated on every row, so they can be
#    formatted as simple blocks; only G has gaps (the 3 "reference" rows)
#    and was already handled cell-by-cell above, together with G1.
# ---------------------------------------------------------------------------
$ws.Range("G1").HorizontalAlignment = -4108  # xlCenter

$alwaysFull = $excel.Union($ws.Range("A1:F24"), $ws.Range("H1:K24"))
$alwaysFull.HorizontalAlignment = -4108  # xlCenter

$ws.Range("J2:K24").WrapText = $true

# ---------------------------------------------------------------------------
# 5) Misc. cosmetic bits that moved in the real edit: selection anchor.
# ---------------------------------------------------------------------------
$ws.Range("G24").Select()
